$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: merge a run of characters at the start of a paragraph (or at
# a known offset) into a single run by re-"finding & replacing" the
# same text within a tightly scoped sub-range. Word's Find/Replace
# engine merges runs that end up with identical formatting, which is
# exactly the effect the diff shows (split runs -> single run).
# ---------------------------------------------------------------------
function Merge-RangeText($startPos, $text) {
    $sub = $d.Range($startPos, $startPos + $text.Length)
    [void]$sub.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2)
}

# 1) Paragraph "14th Apr" heading: runs "1" + "4" -> single run "14"
$p73 = $d.Paragraphs(73).Range
Merge-RangeText $p73.Start "14"

# 2) Paragraph "Research into implementation for first algorithm, Random Forest.":
#    runs "Research" + " into implementation for first algorithm, Random Forest." -> single run
$p74 = $d.Paragraphs(74).Range
Merge-RangeText $p74.Start "Research into implementation for first algorithm, Random Forest."

# 3) Paragraph "19th Apr" heading: runs "1" + "9" -> single run "19"
$p75 = $d.Paragraphs(75).Range
Merge-RangeText $p75.Start "19"

# 4) Paragraph "Further research into implementation ...": four runs -> single run
$p76 = $d.Paragraphs(76).Range
Merge-RangeText $p76.Start "Further research into implementation for first algorithm, Random Forest."

# 5) The trailing "_GoBack" bookmark currently sits at the end of paragraph
#    76; the edit relocates it into the newly written entry below, so
#    remove it from its old spot first (it's recreated further down).
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# 6) Insert the new "20th Apr" journal entry (heading + body + trailing
#    blank paragraph) right after paragraph 76.
$p76 = $d.Paragraphs(76).Range
[void]$p76.InsertParagraphAfter()
$newRange = $d.Paragraphs(77).Range

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$fragment = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document $wns><w:body>
<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>20</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:vertAlign w:val="superscript"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">th </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>Apr</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr><w:t>Initial plan for algorithm written up in human readable language.</w:t></w:r><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr><w:br/><w:t>Email sent to tutor to give update on progress.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Added journal entries for </w:t></w:r><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr><w:t>20</w:t></w:r><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:vertAlign w:val="superscript"/><w:lang w:val="en-US"/></w:rPr><w:t>th</w:t></w:r><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> April.</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

[void]$newRange.InsertXML($fragment)

Write-Output "done"
